$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5733461505019761
$ws.Range("C2").Value = 0.2100473317179308
$ws.Range("D2").Value = 0.07827614943023775
$ws.Range("E2").Value = 0.1232850628349809
$ws.Range("G2").Value = 0.7917072114167638
$ws.Range("H2").Value = 0.8755528824337517
$ws.Range("K2").Value = 0.3051687834335439
$ws.Range("L2").Value = 0.1890576363323362
$ws.Range("O2").Value = 3.350518763409781
$ws.Range("B3").Value = 0.5296011967792538
$ws.Range("C3").Value = 0.2107327780498629
$ws.Range("D3").Value = 0.07104519052593616
$ws.Range("E3").Value = 0.1228729911918407
$ws.Range("G3").Value = 0.7944903458471941
$ws.Range("H3").Value = 0.8810249826384862
$ws.Range("K3").Value = 0.2681359258690179
$ws.Range("L3").Value = 0.1818536565154858
$ws.Range("O3").Value = 3.367582627110309
$ws.Range("B4").Value = 0.5028873827448024
$ws.Range("C4").Value = 0.2111882631674185
$ws.Range("D4").Value = 0.06663876808228508
$ws.Range("E4").Value = 0.1226842674476991
$ws.Range("G4").Value = 0.7966896045556666
$ws.Range("H4").Value = 0.8847546179221126
$ws.Range("K4").Value = 0.2453653395385516
$ws.Range("L4").Value = 0.1775287838650002
$ws.Range("O4").Value = 3.379863279803402
$ws.Range("B5").Value = 0.4920386086767792
$ws.Range("C5").Value = 0.2113826107472399
$ws.Range("D5").Value = 0.06485153824132794
$ws.Range("E5").Value = 0.1226235439022467
$ws.Range("G5").Value = 0.7977090667291336
$ws.Range("H5").Value = 0.8863675090351748
$ws.Range("K5").Value = 0.2360785455962571
$ws.Range("L5").Value = 0.1757911784063708
$ws.Range("O5").Value = 3.385321188037793
$ws.Range("B6").Value = 0.4902394518948654
$ws.Range("C6").Value = 0.2114154103526396
$ws.Range("D6").Value = 0.06455527942486583
$ws.Range("E6").Value = 0.1226144387033514
$ws.Range("G6").Value = 0.797885789391799
$ws.Range("H6").Value = 0.8866409491800766
$ws.Range("K6").Value = 0.2345360370924112
$ws.Range("L6").Value = 0.1755041512422224
$ws.Range("O6").Value = 3.3862548555052
$ws.Range("B7").Value = 0.5027409203800062
$ws.Range("C7").Value = 0.2111908488069112
$ws.Range("D7").Value = 0.06661463074276242
$ws.Range("E7").Value = 0.1226833829623644
$ws.Range("G7").Value = 0.7967028544523203
$ws.Range("H7").Value = 0.884775993165249
$ws.Range("K7").Value = 0.2452401246521134
$ws.Range("L7").Value = 0.1775052493383953
$ws.Range("O7").Value = 3.379935051135391
$ws.Range("B8").Value = 0.558233023199989
$ws.Range("C8").Value = 0.2102765056888494
$ws.Range("D8").Value = 0.07577598949460196
$ws.Range("E8").Value = 0.1231296467372722
$ws.Range("G8").Value = 0.7925650099018782
$ws.Range("H8").Value = 0.8773629575639035
$ws.Range("K8").Value = 0.292406848661301
$ws.Range("L8").Value = 0.1865533215516137
$ws.Range("O8").Value = 3.356028071995212
$ws.Range("B9").Value = 0.6681858198502084
$ws.Range("C9").Value = 0.2087569117783694
$ws.Range("D9").Value = 0.09400657429024761
$ws.Range("E9").Value = 0.1245143477306669
$ws.Range("G9").Value = 0.7883457388871733
$ws.Range("H9").Value = 0.8657573824229985
$ws.Range("K9").Value = 0.3846264504548458
$ws.Range("L9").Value = 0.2050752839892738
$ws.Range("O9").Value = 3.323459104281284
$ws.Range("B10").Value = 0.7496355313140839
$ws.Range("C10").Value = 0.2078054861429735
$ws.Range("D10").Value = 0.1075637239625706
$ws.Range("E10").Value = 0.1258419221257938
$ws.Range("G10").Value = 0.7876268716395174
$ws.Range("H10").Value = 0.8590149352135938
$ws.Range("K10").Value = 0.4521952063437311
$ws.Range("L10").Value = 0.2191570007111352
$ws.Range("O10").Value = 3.308264251711336
$ws.Range("B11").Value = 0.7868296355417783
$ws.Range("C11").Value = 0.20740814273789
$ws.Range("D11").Value = 0.113767035520965
$ws.Range("E11").Value = 0.1265131644392135
$ws.Range("G11").Value = 0.7878183277623663
$ws.Range("H11").Value = 0.8563344649770954
$ws.Range("K11").Value = 0.4828904665495202
$ws.Range("L11").Value = 0.2256658490468908
$ws.Range("O11").Value = 3.303250144203304
$ws.Range("B12").Value = 0.8009339695782671
$ws.Range("C12").Value = 0.2072627513670682
$ws.Range("D12").Value = 0.116121256606263
$ws.Range("E12").Value = 0.1267770152617658
$ws.Range("G12").Value = 0.7879654818219137
$ws.Range("H12").Value = 0.8553749999428391
$ws.Range("K12").Value = 0.494507478693464
$ws.Range("L12").Value = 0.2281453436231544
$ws.Range("O12").Value = 3.301624493867649
$ws.Range("B13").Value = 0.7978954826116933
$ws.Range("C13").Value = 0.2072938387422454
$ws.Range("D13").Value = 0.1156140038771554
$ws.Range("E13").Value = 0.1267197606101824
$ws.Range("G13").Value = 0.7879304677638999
$ws.Range("H13").Value = 0.8555791669543851
$ws.Range("K13").Value = 0.4920058504776534
$ws.Range("L13").Value = 0.2276106859898306
$ws.Range("O13").Value = 3.301962459120233
$ws.Range("B14").Value = 0.7879896161953752
$ws.Range("C14").Value = 0.2073960797596328
$ws.Range("D14").Value = 0.113960615427942
$ws.Range("E14").Value = 0.1265346779713852
$ws.Range("G14").Value = 0.7878289374841785
$ws.Range("H14").Value = 0.8562544156970944
$ws.Range("K14").Value = 0.4838463409499241
$ws.Range("L14").Value = 0.2258695436891287
$ws.Range("O14").Value = 3.303110927190914
$ws.Range("B15").Value = 0.7819245356797353
$ws.Range("C15").Value = 0.2074593653524133
$ws.Range("D15").Value = 0.1129485385855276
$ws.Range("E15").Value = 0.126422567855613
$ws.Range("G15").Value = 0.7877764722758229
$ws.Range("H15").Value = 0.8566752609364698
$ws.Range("K15").Value = 0.4788475274303892
$ws.Range("L15").Value = 0.2248049603937545
$ws.Range("O15").Value = 3.303849964622316
$ws.Range("B16").Value = 0.7472076096310616
$ws.Range("C16").Value = 0.2078321648947536
$ws.Range("D16").Value = 0.1071590468524874
$ws.Range("E16").Value = 0.1257994082825817
$ws.Range("G16").Value = 0.7876248009293221
$ws.Range("H16").Value = 0.8591978885694402
$ws.Range("K16").Value = 0.4501883035358389
$ws.Range("L16").Value = 0.2187336983104018
$ws.Range("O16").Value = 3.308630142020036
$ws.Range("B17").Value = 0.7259458153799869
$ws.Range("C17").Value = 0.2080699302411801
$ws.Range("D17").Value = 0.1036166059284795
$ws.Range("E17").Value = 0.125434352295688
$ws.Range("G17").Value = 0.7876646187995817
$ws.Range("H17").Value = 0.8608444567234557
$ws.Range("K17").Value = 0.4325956130952306
$ws.Range("L17").Value = 0.2150355029821469
$ws.Range("O17").Value = 3.312048875425148
$ws.Range("B18").Value = 0.7137300230383801
$ws.Range("C18").Value = 0.2082100262274409
$ws.Range("D18").Value = 0.1015824869273132
$ws.Range("E18").Value = 0.1252307193175781
$ws.Range("G18").Value = 0.7877363162609754
$ws.Range("H18").Value = 0.8618279196390262
$ws.Range("K18").Value = 0.4224728331202527
$ws.Range("L18").Value = 0.2129180987639785
$ws.Range("O18").Value = 3.314193891131424
$ws.Range("B19").Value = 0.7095962940435641
$ws.Range("C19").Value = 0.2082580347000942
$ws.Range("D19").Value = 0.1008943538909506
$ws.Range("E19").Value = 0.1251628616260305
$ws.Range("G19").Value = 0.7877689694227001
$ws.Range("H19").Value = 0.8621671562171258
$ws.Range("K19").Value = 0.4190447759769711
$ws.Range("L19").Value = 0.2122028503885929
$ws.Range("O19").Value = 3.31495083579938
$ws.Range("B20").Value = 0.7282077849907864
$ws.Range("C20").Value = 0.2080442742629245
$ws.Range("D20").Value = 0.1039933532028812
$ws.Range("E20").Value = 0.1254725573217641
$ws.Range("G20").Value = 0.7876553292764754
$ws.Range("H20").Value = 0.8606654096849127
$ws.Range("K20").Value = 0.4344687958173097
$ws.Range("L20").Value = 0.2154281788403267
$ws.Range("O20").Value = 3.311666454734194
$ws.Range("B21").Value = 0.790898679415875
$ws.Range("C21").Value = 0.2073659115928308
$ws.Range("D21").Value = 0.1144461158436343
$ws.Range("E21").Value = 0.1265887790177622
$ws.Range("G21").Value = 0.7878567325296899
$ws.Range("H21").Value = 0.8560545708853056
$ws.Range("K21").Value = 0.4862431707310861
$ws.Range("L21").Value = 0.2263805601012763
$ws.Range("O21").Value = 3.302766181784847
$ws.Range("B22").Value = 0.8319854962162481
$ws.Range("C22").Value = 0.206952125498816
$ws.Range("D22").Value = 0.1213076724203006
$ws.Range("E22").Value = 0.1273746278396537
$ws.Range("G22").Value = 0.7884235453793735
$ws.Range("H22").Value = 0.8533650139649751
$ws.Range("K22").Value = 0.5200418911666702
$ws.Range("L22").Value = 0.2336244251321205
$ws.Range("O22").Value = 3.298541140501413
$ws.Range("B23").Value = 0.8100464309030713
$ws.Range("C23").Value = 0.207170274350581
$ws.Range("D23").Value = 0.1176427887219944
$ws.Range("E23").Value = 0.1269500559766072
$ws.Range("G23").Value = 0.7880811751652033
$ws.Range("H23").Value = 0.8547708569090702
$ws.Range("K23").Value = 0.5020066157156577
$ws.Range("L23").Value = 0.2297504082053621
$ws.Range("O23").Value = 3.300650429484392
$ws.Range("B24").Value = 0.727185123880389
$ws.Range("C24").Value = 0.2080558627324152
$ws.Range("D24").Value = 0.1038230180422062
$ws.Range("E24").Value = 0.1254552653837315
$ws.Range("G24").Value = 0.7876593770499483
$ws.Range("H24").Value = 0.8607462421266803
$ws.Range("K24").Value = 0.4336219563686541
$ws.Range("L24").Value = 0.215250622864076
$ws.Range("O24").Value = 3.311838787876809
$ws.Range("B25").Value = 0.6383218006850768
$ws.Range("C25").Value = 0.2091389036199018
$ws.Range("D25").Value = 0.08904617297412187
$ws.Range("E25").Value = 0.1240852383855753
$ws.Range("G25").Value = 0.789069497140801
$ws.Range("H25").Value = 0.8685834730420225
$ws.Range("K25").Value = 0.359709811758421
$ws.Range("L25").Value = 0.1999813476374754
$ws.Range("O25").Value = 3.330736758954458

Write-Host "updated 216 cells"
